$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.128355503082275
$ws.Range("C1").Value = 2.89504337310791
$ws.Range("D1").Value = 3.233790636062622
$ws.Range("E1").Value = 15
